$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..13 (A: imageName/id as number, B: centerX, C: centerY)
$data = @(
    @(11000, 196, 126),
    @(11001, 205, 125),
    @(11002, 148, 117),
    @(11003, 168, 88),
    @(11004, 166, 142),
    @(11005, 204, 100),
    @(11006, 135, 114),
    @(11007, 176, 105),
    @(11008, 181, 119),
    @(11009, 158, 104),
    @(11010, 183, 120),
    @(11011, 190, 113)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

$ws.Range("D13").Select()
